# Horarios actualizados Linea 141 - 944
# Scrape timestamp updates: 03:51:22 -> 04:15:01

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "LP1912": the bus that already arrived (81_EL PELIGRO @04:02)
# dropped off the list, every remaining row shifted up one position,
# and four new arrivals were appended at the bottom.
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Remove the row for the bus that already departed; this shifts the
# rest of the data rows up by one, just like the source scrape does.
$ws1.Rows.Item(6).Delete()

$ws1.Range("A2").Value = "Última actualización: 04:15:01"
$ws1.Range("A3").Value = "Total filas: 10"

$data1 = @(
    @("04:15:01", "04:47", "215_EL PELIGRO", 32, "LP1912"),
    @("04:15:01", "04:53", "11_ETCHEVERRY", 38, "LP1912"),
    @("04:15:01", "05:12", "17_ROMERO", 57, "LP1912"),
    @("04:15:01", "05:22", "23_HERNANDEZ", 67, "LP1912"),
    @("04:15:01", "05:32", "81_EL PELIGRO", 77, "LP1912"),
    @("04:15:01", "05:44", "14_ABASTO", 89, "LP1912"),
    @("04:15:01", "05:52", "17_ROMERO", 97, "LP1912"),
    @("04:15:01", "06:01", "16_SANTA ANA", 106, "LP1912"),
    @("04:15:01", "06:04", "10_OLMOS", 109, "LP1912"),
    @("04:15:01", "06:11", "215A_EL PATO", 116, "LP1912")
)

$row = 6
foreach ($d in $data1) {
    $ws1.Range("A$row").Value = $d[0]
    $ws1.Range("B$row").Value = $d[1]
    $ws1.Range("C$row").Value = $d[2]
    $ws1.Range("D$row").Value = $d[3]
    $ws1.Range("E$row").Value = $d[4]
    $row = $row + 1
}

# -----------------------------------------------------------------
# Sheet "LP1912-215": only line 215 buses. The existing arrival stays
# and a second (215A) arrival is appended.
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:15:01"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = "04:15:01"
$ws2.Range("D6").Value = 32

$ws2.Range("A7").Value = "04:15:01"
$ws2.Range("B7").Value = "06:11"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 116
$ws2.Range("E7").Value = "LP1912"

# -----------------------------------------------------------------
# Sheet "6203-6173": no arrivals, only the refresh timestamp changes.
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 04:15:01"
